# chore: status log + auto-updated Current Price
# Append a new status-log entry (row 3) to the Status_Log sheet, mirroring
# the structure of the existing row 2 entry with a refreshed timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Status_Log")

$row = 3

$ws.Cells.Item($row, 1).Value = "2025-12-29T06:33:33Z"
$ws.Cells.Item($row, 2).Value = 80.40000152587891
$ws.Cells.Item($row, 3).Value = 15.31404656443955
$ws.Cells.Item($row, 4).Value = 425.008208559155
$ws.Cells.Item($row, 5).Value = 0
$ws.Cells.Item($row, 6).Value = 0
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = 0
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = "观察"
$ws.Cells.Item($row, 11).Value = 0
